$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve trailing zeros / exact formatting)
$textFormatCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D13', 'D14', 'D16', 'D17', 'D22', 'D23', 'D25', 'D26', 'D31', 'D32', 'D33', 'D35', 'D36', 'D38', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D49', 'D51')
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values row by row
$ws.Range("D2").Value = '42.924.45'
$ws.Range("E2").Value = '  +4.21%  '

$ws.Range("D3").Value = '2.285.63'
$ws.Range("E3").Value = '  +5.11%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '252.25'
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").Value = '0.641'
$ws.Range("E6").Value = '  +4.65%  '

$ws.Range("D7").Value = '72.85'
$ws.Range("E7").Value = '  +10.27%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '0.661'
$ws.Range("E9").Value = '  +14.86%  '

$ws.Range("D10").Value = '38.99'
$ws.Range("E10").Value = '  +7.71%  '

$ws.Range("D12").Value = '59.94'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").Value = '7.38'
$ws.Range("E13").Value = '  +8.21%  '

$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  +2.42%  '

$ws.Range("D15").Value = '2.625.62'
$ws.Range("E15").Value = '  +5.11%  '

$ws.Range("D16").Value = '15.08'
$ws.Range("E16").Value = '  +5.84%  '

$ws.Range("D17").Value = '0.893'
$ws.Range("E17").Value = '  +5.77%  '

$ws.Range("D18").Value = '2.281.03'
$ws.Range("E18").Value = '  +4.78%  '

$ws.Range("D19").Value = '42.859.74'
$ws.Range("E19").Value = '  +4.27%  '

$ws.Range("E20").Value = '  +7.44%  '

$ws.Range("E21").Value = '  +5.68%  '

$ws.Range("D22").Value = '73.53'
$ws.Range("E22").Value = '  +2.83%  '

$ws.Range("D23").Value = '237.78'
$ws.Range("E23").Value = '  +3.49%  '

$ws.Range("E24").Value = '  +7.06%  '

$ws.Range("D25").Value = '3.88'
$ws.Range("E25").Value = '  +1.80%  '

$ws.Range("D26").Value = '11.67'
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("E27").Value = '  -0.12%  '

$ws.Range("E28").Value = '  +1.91%  '

$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("E30").Value = '  +3.78%  '

$ws.Range("D31").Value = '168.01'
$ws.Range("E31").Value = '  -0.30%  '

$ws.Range("D32").Value = '21.12'
$ws.Range("E32").Value = '  +4.78%  '

$ws.Range("D33").Value = '6.31'
$ws.Range("E33").Value = '  +11.18%  '

$ws.Range("E34").Value = '  +6.32%  '

$ws.Range("D35").Value = '0.0815'
$ws.Range("E35").Value = '  +9.04%  '

$ws.Range("D36").Value = '31.08'
$ws.Range("E36").Value = '  +27.40%  '

$ws.Range("E37").Value = '  +5.10%  '

$ws.Range("D38").Value = '4.77'
$ws.Range("E38").Value = '  +21.65%  '

$ws.Range("E39").Value = '  +6.35%  '

$ws.Range("E40").Value = '  +1.48%  '

$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '2.33'
$ws.Range("E41").Value = '  +5.79%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '13.33'
$ws.Range("E42").Value = '  +17.96%  '

$ws.Range("D43").Value = '6.04'
$ws.Range("E43").Value = '  +10.31%  '

$ws.Range("E44").Value = '  +14.75%  '

$ws.Range("D45").Value = '9.23'
$ws.Range("E45").Value = '  +9.02%  '

$ws.Range("D46").Value = '5.02'
$ws.Range("E46").Value = '  -7.50%  '

$ws.Range("D47").Value = '61.71'
$ws.Range("E47").Value = '  +1.83%  '

$ws.Range("E48").Value = '  +5.35%  '

$ws.Range("D49").Value = '1.20'
$ws.Range("E49").Value = '  +5.11%  '

$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '97.98'
$ws.Range("E51").Value = '  +8.39%  '

